$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 144
$rowCount = 44
$endRow = $startRow + $rowCount - 1

# Columns A (Any), B (Fotograma), J (Projecte), K (Rollo) -- no new shared strings introduced here
$dataABJK = New-Object 'object[,]' $rowCount,11
$dataABJK[0,0] = 1814
$dataABJK[0,1] = 3
$dataABJK[0,9] = 'SPN 2,02 C'
$dataABJK[0,10] = 47
$dataABJK[1,0] = 1814
$dataABJK[1,1] = 3
$dataABJK[1,9] = 'SPN 2,02 C'
$dataABJK[1,10] = 47
$dataABJK[2,0] = 1814
$dataABJK[2,1] = 4
$dataABJK[2,9] = 'SPN 2,02 C'
$dataABJK[2,10] = 47
$dataABJK[3,0] = 1814
$dataABJK[3,1] = 4
$dataABJK[3,9] = 'SPN 2,02 C'
$dataABJK[3,10] = 47
$dataABJK[4,0] = 1814
$dataABJK[4,1] = 4
$dataABJK[4,9] = 'SPN 2,02 C'
$dataABJK[4,10] = 47
$dataABJK[5,0] = 1814
$dataABJK[5,1] = 5
$dataABJK[5,9] = 'SPN 2,02 C'
$dataABJK[5,10] = 47
$dataABJK[6,0] = 1814
$dataABJK[6,1] = 5
$dataABJK[6,9] = 'SPN 2,02 C'
$dataABJK[6,10] = 47
$dataABJK[7,0] = 1814
$dataABJK[7,1] = 5
$dataABJK[7,9] = 'SPN 2,02 C'
$dataABJK[7,10] = 47
$dataABJK[8,0] = 1814
$dataABJK[8,1] = 6
$dataABJK[8,9] = 'SPN 2,02 C'
$dataABJK[8,10] = 47
$dataABJK[9,0] = 1815
$dataABJK[9,1] = 6
$dataABJK[9,9] = 'SPN 2,02 C'
$dataABJK[9,10] = 47
$dataABJK[10,0] = 1815
$dataABJK[10,1] = 6
$dataABJK[10,9] = 'SPN 2,02 C'
$dataABJK[10,10] = 47
$dataABJK[11,0] = 1815
$dataABJK[11,1] = 7
$dataABJK[11,9] = 'SPN 2,02 C'
$dataABJK[11,10] = 47
$dataABJK[12,0] = 1816
$dataABJK[12,1] = 7
$dataABJK[12,9] = 'SPN 2,02 C'
$dataABJK[12,10] = 47
$dataABJK[13,0] = 1816
$dataABJK[13,1] = 8
$dataABJK[13,9] = 'SPN 2,02 C'
$dataABJK[13,10] = 47
$dataABJK[14,0] = 1816
$dataABJK[14,1] = 8
$dataABJK[14,9] = 'SPN 2,02 C'
$dataABJK[14,10] = 47
$dataABJK[15,0] = 1816
$dataABJK[15,1] = 8
$dataABJK[15,9] = 'SPN 2,02 C'
$dataABJK[15,10] = 47
$dataABJK[16,0] = 1816
$dataABJK[16,1] = 9
$dataABJK[16,9] = 'SPN 2,02 C'
$dataABJK[16,10] = 47
$dataABJK[17,0] = 1816
$dataABJK[17,1] = 9
$dataABJK[17,9] = 'SPN 2,02 C'
$dataABJK[17,10] = 47
$dataABJK[18,0] = 1816
$dataABJK[18,1] = 10
$dataABJK[18,9] = 'SPN 2,02 C'
$dataABJK[18,10] = 47
$dataABJK[19,0] = 1816
$dataABJK[19,1] = 10
$dataABJK[19,9] = 'SPN 2,02 C'
$dataABJK[19,10] = 47
$dataABJK[20,0] = 1816
$dataABJK[20,1] = 11
$dataABJK[20,9] = 'SPN 2,02 C'
$dataABJK[20,10] = 47
$dataABJK[21,0] = 1817
$dataABJK[21,1] = 11
$dataABJK[21,9] = 'SPN 2,02 C'
$dataABJK[21,10] = 47
$dataABJK[22,0] = 1817
$dataABJK[22,1] = 11
$dataABJK[22,9] = 'SPN 2,02 C'
$dataABJK[22,10] = 47
$dataABJK[23,0] = 1817
$dataABJK[23,1] = 12
$dataABJK[23,9] = 'SPN 2,02 C'
$dataABJK[23,10] = 47
$dataABJK[24,0] = 1817
$dataABJK[24,1] = 12
$dataABJK[24,9] = 'SPN 2,02 C'
$dataABJK[24,10] = 47
$dataABJK[25,0] = 1817
$dataABJK[25,1] = 13
$dataABJK[25,9] = 'SPN 2,02 C'
$dataABJK[25,10] = 47
$dataABJK[26,0] = 1817
$dataABJK[26,1] = 13
$dataABJK[26,9] = 'SPN 2,02 C'
$dataABJK[26,10] = 47
$dataABJK[27,0] = 1817
$dataABJK[27,1] = 14
$dataABJK[27,9] = 'SPN 2,02 C'
$dataABJK[27,10] = 47
$dataABJK[28,0] = 1818
$dataABJK[28,1] = 14
$dataABJK[28,9] = 'SPN 2,02 C'
$dataABJK[28,10] = 47
$dataABJK[29,0] = 1818
$dataABJK[29,1] = 14
$dataABJK[29,9] = 'SPN 2,02 C'
$dataABJK[29,10] = 47
$dataABJK[30,0] = 1818
$dataABJK[30,1] = 15
$dataABJK[30,9] = 'SPN 2,02 C'
$dataABJK[30,10] = 47
$dataABJK[31,0] = 1818
$dataABJK[31,1] = 15
$dataABJK[31,9] = 'SPN 2,02 C'
$dataABJK[31,10] = 47
$dataABJK[32,0] = 1819
$dataABJK[32,1] = 16
$dataABJK[32,9] = 'SPN 2,02 C'
$dataABJK[32,10] = 47
$dataABJK[33,0] = 1819
$dataABJK[33,1] = 16
$dataABJK[33,9] = 'SPN 2,02 C'
$dataABJK[33,10] = 47
$dataABJK[34,0] = 1819
$dataABJK[34,1] = 17
$dataABJK[34,9] = 'SPN 2,02 C'
$dataABJK[34,10] = 47
$dataABJK[35,0] = 1819
$dataABJK[35,1] = 17
$dataABJK[35,9] = 'SPN 2,02 C'
$dataABJK[35,10] = 47
$dataABJK[36,0] = 1819
$dataABJK[36,1] = 17
$dataABJK[36,9] = 'SPN 2,02 C'
$dataABJK[36,10] = 47
$dataABJK[37,0] = 1819
$dataABJK[37,1] = 18
$dataABJK[37,9] = 'SPN 2,02 C'
$dataABJK[37,10] = 47
$dataABJK[38,0] = 1819
$dataABJK[38,1] = 18
$dataABJK[38,9] = 'SPN 2,02 C'
$dataABJK[38,10] = 47
$dataABJK[39,0] = 1820
$dataABJK[39,1] = 19
$dataABJK[39,9] = 'SPN 2,02 C'
$dataABJK[39,10] = 47
$dataABJK[40,0] = 1820
$dataABJK[40,1] = 19
$dataABJK[40,9] = 'SPN 2,02 C'
$dataABJK[40,10] = 47
$dataABJK[41,0] = 1820
$dataABJK[41,1] = 19
$dataABJK[41,9] = 'SPN 2,02 C'
$dataABJK[41,10] = 47
$dataABJK[42,0] = 1820
$dataABJK[42,1] = 20
$dataABJK[42,9] = 'SPN 2,02 C'
$dataABJK[42,10] = 47
$dataABJK[43,0] = 1820
$dataABJK[43,1] = 20
$dataABJK[43,9] = 'SPN 2,02 C'
$dataABJK[43,10] = 47
$ws.Range("A${startRow}:K${endRow}").Value = $dataABJK

# Column M (Anys) filled next so its shared string is registered before L and C
$dataM = New-Object 'object[,]' $rowCount,1
$dataM[0,0] = '1814-1851'
$dataM[1,0] = '1814-1851'
$dataM[2,0] = '1814-1851'
$dataM[3,0] = '1814-1851'
$dataM[4,0] = '1814-1851'
$dataM[5,0] = '1814-1851'
$dataM[6,0] = '1814-1851'
$dataM[7,0] = '1814-1851'
$dataM[8,0] = '1814-1851'
$dataM[9,0] = '1814-1851'
$dataM[10,0] = '1814-1851'
$dataM[11,0] = '1814-1851'
$dataM[12,0] = '1814-1851'
$dataM[13,0] = '1814-1851'
$dataM[14,0] = '1814-1851'
$dataM[15,0] = '1814-1851'
$dataM[16,0] = '1814-1851'
$dataM[17,0] = '1814-1851'
$dataM[18,0] = '1814-1851'
$dataM[19,0] = '1814-1851'
$dataM[20,0] = '1814-1851'
$dataM[21,0] = '1814-1851'
$dataM[22,0] = '1814-1851'
$dataM[23,0] = '1814-1851'
$dataM[24,0] = '1814-1851'
$dataM[25,0] = '1814-1851'
$dataM[26,0] = '1814-1851'
$dataM[27,0] = '1814-1851'
$dataM[28,0] = '1814-1851'
$dataM[29,0] = '1814-1851'
$dataM[30,0] = '1814-1851'
$dataM[31,0] = '1814-1851'
$dataM[32,0] = '1814-1851'
$dataM[33,0] = '1814-1851'
$dataM[34,0] = '1814-1851'
$dataM[35,0] = '1814-1851'
$dataM[36,0] = '1814-1851'
$dataM[37,0] = '1814-1851'
$dataM[38,0] = '1814-1851'
$dataM[39,0] = '1814-1851'
$dataM[40,0] = '1814-1851'
$dataM[41,0] = '1814-1851'
$dataM[42,0] = '1814-1851'
$dataM[43,0] = '1814-1851'
$ws.Range("M${startRow}:M${endRow}").Value = $dataM

# Column L (Serie) filled next
$dataL = New-Object 'object[,]' $rowCount,1
$dataL[0,0] = 'A,4'
$dataL[1,0] = 'A,4'
$dataL[2,0] = 'A,4'
$dataL[3,0] = 'A,4'
$dataL[4,0] = 'A,4'
$dataL[5,0] = 'A,4'
$dataL[6,0] = 'A,4'
$dataL[7,0] = 'A,4'
$dataL[8,0] = 'A,4'
$dataL[9,0] = 'A,4'
$dataL[10,0] = 'A,4'
$dataL[11,0] = 'A,4'
$dataL[12,0] = 'A,4'
$dataL[13,0] = 'A,4'
$dataL[14,0] = 'A,4'
$dataL[15,0] = 'A,4'
$dataL[16,0] = 'A,4'
$dataL[17,0] = 'A,4'
$dataL[18,0] = 'A,4'
$dataL[19,0] = 'A,4'
$dataL[20,0] = 'A,4'
$dataL[21,0] = 'A,4'
$dataL[22,0] = 'A,4'
$dataL[23,0] = 'A,4'
$dataL[24,0] = 'A,4'
$dataL[25,0] = 'A,4'
$dataL[26,0] = 'A,4'
$dataL[27,0] = 'A,4'
$dataL[28,0] = 'A,4'
$dataL[29,0] = 'A,4'
$dataL[30,0] = 'A,4'
$dataL[31,0] = 'A,4'
$dataL[32,0] = 'A,4'
$dataL[33,0] = 'A,4'
$dataL[34,0] = 'A,4'
$dataL[35,0] = 'A,4'
$dataL[36,0] = 'A,4'
$dataL[37,0] = 'A,4'
$dataL[38,0] = 'A,4'
$dataL[39,0] = 'A,4'
$dataL[40,0] = 'A,4'
$dataL[41,0] = 'A,4'
$dataL[42,0] = 'A,4'
$dataL[43,0] = 'A,4'
$ws.Range("L${startRow}:L${endRow}").Value = $dataL

# Column C (Cognoms Familia) filled last
$dataC = New-Object 'object[,]' $rowCount,1
$dataC[0,0] = 'Margall Busquets'
$dataC[1,0] = 'Vilamajor Civit'
$dataC[2,0] = 'Vilamajor Codina'
$dataC[3,0] = 'Vinaixa Cascalló'
$dataC[4,0] = 'Batlle Asqueró'
$dataC[5,0] = 'Bellet Martí'
$dataC[6,0] = 'Solsona Jovells'
$dataC[7,0] = 'Vilaplana Civit'
$dataC[8,0] = 'Ariatós Giné'
$dataC[9,0] = 'Roiger Caelles'
$dataC[10,0] = 'Lamarca Vergé'
$dataC[11,0] = 'Bertran Farré'
$dataC[12,0] = 'Ponsarnau Torrent'
$dataC[13,0] = 'Novell Soler'
$dataC[14,0] = 'Civit Carrera'
$dataC[15,0] = 'Boldú Violant'
$dataC[16,0] = 'Soler Giné'
$dataC[17,0] = 'Bresolí Cascalló'
$dataC[18,0] = 'Trepat Giner'
$dataC[19,0] = 'Giné Pedrós'
$dataC[20,0] = 'Novell Torrent'
$dataC[21,0] = 'Mas Vergé'
$dataC[22,0] = 'Berniell Gené'
$dataC[23,0] = 'Niubó Martí'
$dataC[24,0] = 'Llas Espervé'
$dataC[25,0] = 'Sucarrat Mirassó'
$dataC[26,0] = 'Torrent Cascalló'
$dataC[27,0] = 'Font Coll'
$dataC[28,0] = 'Capdevila Mas'
$dataC[29,0] = 'Soler Mas'
$dataC[30,0] = 'Rossell Agulló'
$dataC[31,0] = 'Ginestà Mosset'
$dataC[32,0] = 'Oromi Roma'
$dataC[33,0] = 'Mas Viladebait'
$dataC[34,0] = 'Martí Batlle'
$dataC[35,0] = 'Cisteré Mosset'
$dataC[36,0] = 'Planes Pallaas'
$dataC[37,0] = 'Roig Coll'
$dataC[38,0] = 'Monyart Bellera'
$dataC[39,0] = 'Agulló Calderó'
$dataC[40,0] = 'Civit Mas'
$dataC[41,0] = 'Batlle Pujol'
$dataC[42,0] = 'Caelles Mata'
$dataC[43,0] = 'Mas Bobé'
$ws.Range("C${startRow}:C${endRow}").Value = $dataC

# Update selection / active cell to the row after the last new entry
$selAddr = "A" + ($endRow + 1)
$ws.Range($selAddr).Select() | Out-Null

# Scroll the view so the newly added rows are visible (best effort)
$win = $excel.ActiveWindow
$win.ScrollRow = $startRow + 6
$win.ScrollColumn = 1
